$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to "BetaFiberA"
$ws.Name = "BetaFiberA"

# Append a new row (16) mirroring the layout/formatting of row 15, holding
# the averaged intensities computed with the new Gaussian Quadrature Scheme.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)   # xlPasteFormats - copy style s="1" only
$ws.Range("A16").Value2 = 14

$ws.Range("B16").Value2 = $ws.Range("B15").Text

$ws.Range("C16").Value2 = 1.006681252277507
$ws.Range("D16").Value2 = 0.9668042872944713
$ws.Range("E16").Value2 = 0.9982645422461321
$ws.Range("F16").Value2 = 0.9918008236398924
$ws.Range("G16").Value2 = 1.006681252277507
$ws.Range("H16").Value2 = 0.9668042872944713
$ws.Range("I16").Value2 = 1.002759702876757
$ws.Range("J16").Value2 = 0.9879085223676277
$ws.Range("K16").Value2 = 0.998656679882176
$ws.Range("L16").Value2 = 0.9761799342428142
$ws.Range("M16").Value2 = 1.006681252277507
$ws.Range("N16").Value2 = 0.9825344147703017
$ws.Range("O16").Value2 = 0.9908877263645007
$ws.Range("P16").Value2 = 0.9911319681034223
